$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: add D7 = "true" (column D = integration).
# The existing C7 cell already holds the text value "true" (not a boolean),
# so copy/paste it across to keep the same text-typed cell instead of
# letting a fresh "true" entry get auto-converted into a boolean.
$ws.Range("C7").Copy()
$ws.Range("D7").PasteSpecial()

# Row 8: new test case row
$ws.Range("B8").Value = "Text area test"
$ws.Range("C7").Copy()
$ws.Range("C8").PasteSpecial()
$ws.Range("D7").Copy()
$ws.Range("D8").PasteSpecial()
$ws.Range("E8").Value = "такое себе"

$excel.CutCopyMode = $false

# Update selection to E9 as in the diff
$ws.Range("E9").Select()
